# ---------------------------------------------------------------------------
# Applies the "cambios en los analisis" commit:
#   1. Swaps the LSPM / LSPMW labels (they were transposed) wherever they
#      appear as row/column headers: Matriz_Resultados, P_valores,
#      Estadisticos_HLN_DM (col C/D headers + row A3/A4) and Resumen_Modelos
#      (row A3/A4).
#   2. Rewrites the P_valores (p-value) and Estadisticos_HLN_DM (HLN
#      statistic) 9x9 matrices with the recomputed values that result from
#      that swap (LSPM and LSPMW were fit again / re-paired, so the whole
#      block of numbers shifts, not just the two affected rows/columns).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix the transposed LSPM / LSPMW labels.
# ---------------------------------------------------------------------------
$matrixSheets = @("Matriz_Resultados", "P_valores", "Estadisticos_HLN_DM")
foreach ($name in $matrixSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C1").Value2 = "LSPMW"
    $ws.Range("D1").Value2 = "LSPM"
    $ws.Range("A3").Value2 = "LSPMW"
    $ws.Range("A4").Value2 = "LSPM"
}

$wsResumen = $wb.Worksheets.Item("Resumen_Modelos")
$wsResumen.Range("A3").Value2 = "LSPMW"
$wsResumen.Range("A4").Value2 = "LSPM"

# ---------------------------------------------------------------------------
# 2. Rewrite the two 9x9 result matrices (rows 2-10, columns B-J).
# ---------------------------------------------------------------------------
$PValoresData = @(
    @([double]"1",[double]"2.565412327015792E-10",[double]"1.176800212832063E-09",[double]"9.887793694929314E-09",[double]"3.581022633980524E-08",[double]"7.692991593755139E-05",[double]"4.849504842230346E-05",[double]"0.0008448892605139591",[double]"0.001583936941778274"),
    @([double]"2.565412327015792E-10",[double]"1",[double]"0.04857386869698344",[double]"0.5268464423914816",[double]"0.1800657254820428",[double]"9.383924748362915E-10",[double]"0.2558204081525433",[double]"0.1090030560728956",[double]"0.1218174779132404"),
    @([double]"1.176800212832063E-09",[double]"0.04857386869698344",[double]"1",[double]"0.1126513462053611",[double]"0.028017883212343",[double]"5.290516691403013E-09",[double]"0.5409572751659724",[double]"0.1658459731712094",[double]"0.1767213519189381"),
    @([double]"9.887793694929314E-09",[double]"0.5268464423914816",[double]"0.1126513462053611",[double]"1",[double]"0.2574094495959982",[double]"1.090636803091627E-08",[double]"0.06606342043444902",[double]"0.05814683094478501",[double]"0.06536442873693726"),
    @([double]"3.581022633980524E-08",[double]"0.1800657254820428",[double]"0.028017883212343",[double]"0.2574094495959982",[double]"1",[double]"1.370768452080284E-07",[double]"0.05914687881925218",[double]"0.04247036289500117",[double]"0.04954516074775639"),
    @([double]"7.692991593755139E-05",[double]"9.383924748362915E-10",[double]"5.290516691403013E-09",[double]"1.090636803091627E-08",[double]"1.370768452080284E-07",[double]"1",[double]"7.00116645138138E-05",[double]"0.001247846599785829",[double]"0.002278097591139838"),
    @([double]"4.849504842230346E-05",[double]"0.2558204081525433",[double]"0.5409572751659724",[double]"0.06606342043444902",[double]"0.05914687881925218",[double]"7.00116645138138E-05",[double]"1",[double]"0.1814141569037915",[double]"0.1806883130524053"),
    @([double]"0.0008448892605139591",[double]"0.1090030560728956",[double]"0.1658459731712094",[double]"0.05814683094478501",[double]"0.04247036289500117",[double]"0.001247846599785829",[double]"0.1814141569037915",[double]"1",[double]"0.6100945777787632"),
    @([double]"0.001583936941778274",[double]"0.1218174779132404",[double]"0.1767213519189381",[double]"0.06536442873693726",[double]"0.04954516074775639",[double]"0.002278097591139838",[double]"0.1806883130524053",[double]"0.6100945777787632",[double]"1")
)

$EstadisticosData = @(
    @([double]"0",[double]"-10.87575213492686",[double]"-10.01400570283196",[double]"-8.886130496577135",[double]"-8.242775358544435",[double]"-4.843940959691864",[double]"-5.034232164166514",[double]"-3.861538081875873",[double]"-3.601971362408543"),
    @([double]"10.87575213492686",[double]"0",[double]"-2.088148253900181",[double]"0.6430366127503121",[double]"1.384575484511581",[double]"10.13905072359772",[double]"-1.166706719855327",[double]"-1.670424102934224",[double]"-1.609255537431357"),
    @([double]"10.01400570283196",[double]"2.088148253900181",[double]"0",[double]"1.652418434545484",[double]"2.352287392832455",[double]"9.208891216435564",[double]"-0.6210365308978103",[double]"-1.433249207042692",[double]"-1.395739380823915"),
    @([double]"8.886130496577135",[double]"-0.6430366127503121",[double]"-1.652418434545484",[double]"0",[double]"1.162698386737353",[double]"8.836147043517117",[double]"-1.934124709286154",[double]"-1.998693310036997",[double]"-1.939542201380715"),
    @([double]"8.242775358544435",[double]"-1.384575484511581",[double]"-2.352287392832455",[double]"-1.162698386737353",[double]"0",[double]"7.599203697563123",[double]"-1.990121930914116",[double]"-2.153847724193407",[double]"-2.078385791820676"),
    @([double]"4.843940959691864",[double]"-10.13905072359772",[double]"-9.208891216435564",[double]"-8.836147043517117",[double]"-7.599203697563123",[double]"0",[double]"-4.8827344270493",[double]"-3.700737906976864",[double]"-3.450663116837194"),
    @([double]"5.034232164166514",[double]"1.166706719855327",[double]"0.6210365308978103",[double]"1.934124709286154",[double]"1.990121930914116",[double]"4.8827344270493",[double]"0",[double]"-1.380120857081358",[double]"-1.382515448081809"),
    @([double]"3.861538081875873",[double]"1.670424102934224",[double]"1.433249207042692",[double]"1.998693310036997",[double]"2.153847724193407",[double]"3.700737906976864",[double]"1.380120857081358",[double]"0",[double]"-0.5173211836388842"),
    @([double]"3.601971362408543",[double]"1.609255537431357",[double]"1.395739380823915",[double]"1.939542201380715",[double]"2.078385791820676",[double]"3.450663116837194",[double]"1.382515448081809",[double]"0.5173211836388842",[double]"0")
)

$wsP = $wb.Worksheets.Item("P_valores")
for ($r = 0; $r -lt 9; $r++) {
    for ($c = 0; $c -lt 9; $c++) {
        $wsP.Cells.Item($r + 2, $c + 2).Value2 = $PValoresData[$r][$c]
    }
}

$wsE = $wb.Worksheets.Item("Estadisticos_HLN_DM")
for ($r = 0; $r -lt 9; $r++) {
    for ($c = 0; $c -lt 9; $c++) {
        $wsE.Cells.Item($r + 2, $c + 2).Value2 = $EstadisticosData[$r][$c]
    }
}
